$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "IT - OPEN ACTION ITEM LOG"
$ws.Range("A2").Value = "Project: IT Requirements Traceability Matrix"
$ws.Range("E12").Value = "Ethics Committee"
$ws.Range("E13").Value = "Compliance Officers"
